$d = $word.ActiveDocument

$replacements = @(
    @{old = "999÷4="; new = "269÷5="},
    @{old = "829÷5="; new = "748÷4="},
    @{old = "793÷9="; new = "823÷4="},
    @{old = "620÷7="; new = "985÷7="},
    @{old = "249÷6="; new = "882÷5="},
    @{old = "132÷5="; new = "624÷2="},
    @{old = "623÷8="; new = "595÷2="},
    @{old = "528÷9="; new = "779÷3="},
    @{old = "465÷5="; new = "889÷7="},
    @{old = "335÷8="; new = "979÷3="},
    @{old = "313÷6="; new = "855÷9="},
    @{old = "956÷7="; new = "104÷7="},
    @{old = "442÷4="; new = "388÷2="},
    @{old = "227÷9="; new = "391÷2="},
    @{old = "471÷6="; new = "395÷9="},
    @{old = "267÷9="; new = "946÷2="},
    @{old = "858÷9="; new = "825÷2="},
    @{old = "583÷3="; new = "690÷9="},
    @{old = "564÷2="; new = "995÷7="},
    @{old = "824÷6="; new = "855÷9="},
    @{old = "390÷6="; new = "503÷6="},
    @{old = "290÷9="; new = "202÷2="},
    @{old = "580÷2="; new = "546÷8="},
    @{old = "556÷3="; new = "329÷9="},
    @{old = "303÷4="; new = "591÷8="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
